# Adds proofing marks (grammar/spelling "squiggle" markers) around the
# last word/phrase of a number of notes paragraphs, and appends a new
# "Different perspectives" block of notes at the end of the document.
#
# Word marks runs that the grammar/spell checker flagged by wrapping
# them in <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>
# (or spellStart/spellEnd) and splitting the paragraph's text run so the
# flagged word/phrase is its own run. We reproduce that exact structure
# by deleting each target paragraph's text and re-inserting the
# equivalent OOXML via Range.InsertXML.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Replace-ParagraphXml {
    param(
        [string]$FindText,
        [string]$InnerXml
    )
    $r = $d.Content
    $found = $r.Find.Execute($FindText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $FindText"
    }
    $r.Delete()
    $r.InsertXML("<w:p $wNs>$InnerXml</w:p>")
}

# --- "Engineering:" section -------------------------------------------------

Replace-ParagraphXml `
    "to designing and building things" `
    '<w:r><w:t xml:space="preserve">to designing and building </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>things</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "involves trade-offs, risk analysis, skill, and knowledge" `
    '<w:r><w:t xml:space="preserve">involves trade-offs, risk analysis, skill, and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>knowledge</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# --- "A discipline whose aim ..." paragraph --------------------------------

Replace-ParagraphXml `
    "A discipline whose aim is the production of fault-free software, delivered on-time and within budget, that satisfies the user’s needs. Furthermore, the software must be easy to modify when the user’s needs change." `
    '<w:r><w:t xml:space="preserve">A discipline whose aim is the production of fault-free software, delivered </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>on-time</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and within budget, that satisfies the user’s needs. Furthermore, the software must be easy to modify when the user’s needs change.</w:t></w:r>'

# --- "Buy vs build" ----------------------------------------------------------

Replace-ParagraphXml `
    "Buy vs build" `
    '<w:r><w:t xml:space="preserve">Buy vs </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>build</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# --- "Software engineering can be looked at different perspectives" --------

Replace-ParagraphXml `
    "Software engineering can be looked at different perspectives" `
    '<w:r><w:t xml:space="preserve">Software engineering can be looked at different </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>perspectives</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# --- "Principles of software engineering" list ------------------------------

Replace-ParagraphXml `
    "make quality number 1" `
    '<w:r><w:t xml:space="preserve">make quality number </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>1</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "give products to customers early" `
    '<w:r><w:t xml:space="preserve">give products to customers </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>early</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "understand the problem first" `
    '<w:r><w:t xml:space="preserve">understand the problem </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>first</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "choose an appropriate process" `
    '<w:r><w:t xml:space="preserve">choose an appropriate </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>process</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "good management is more important than good technology" `
    '<w:r><w:t xml:space="preserve">good management is more important than good </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>technology</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "people are the key to success" `
    '<w:r><w:t xml:space="preserve">people are the key to </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>success</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "follow an architecture first process" `
    '<w:r><w:t xml:space="preserve">follow an architecture first </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>process</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "use component based development to reduce coding effort" `
    '<w:r><w:t xml:space="preserve">use </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>component based</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> development to reduce coding effort</w:t></w:r>'

Replace-ParagraphXml `
    "show the customer preliminary versions of the software frequently" `
    '<w:r><w:t xml:space="preserve">show the customer preliminary versions of the software </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>frequently</w:t></w:r><w:proofErr w:type="gramEnd"/>'

Replace-ParagraphXml `
    "have incremental releases" `
    '<w:r><w:t xml:space="preserve">have incremental </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>releases</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# --- "Divide and conquer" (has a leading tab run) ---------------------------

Replace-ParagraphXml `
    "Divide and conquer" `
    '<w:r><w:tab/><w:t xml:space="preserve">Divide and </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>conquer</w:t></w:r><w:proofErr w:type="gramEnd"/>'

# --- Append new "Different perspectives" notes after "Incrementality" ------

$r = $d.Content
$found = $r.Find.Execute("Incrementality", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find text: Incrementality"
}
$insertionPoint = $d.Range($r.End, $r.End)

$newBlockXml = ""
$newBlockXml += "<w:p $wNs/>"
$newBlockXml += "<w:p $wNs><w:r><w:t xml:space=`"preserve`">Different </w:t></w:r><w:proofErr w:type=`"spellStart`"/><w:r><w:t>perspertives</w:t></w:r><w:proofErr w:type=`"spellEnd`"/></w:p>"
$newBlockXml += "<w:p $wNs><w:r><w:tab/><w:t>Business</w:t></w:r></w:p>"
$newBlockXml += "<w:p $wNs><w:r><w:tab/><w:t>Engineering</w:t></w:r></w:p>"
$newBlockXml += "<w:p $wNs><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>User</w:t></w:r></w:p>"
$newBlockXml += "<w:p $wNs><w:r><w:tab/><w:t>Others</w:t></w:r></w:p>"
$newBlockXml += "<w:p $wNs/>"

$insertionPoint.InsertXML($newBlockXml)
